$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("B2").Value = 0.1736111111111111
$ws.Range("C2").Value = 0.5902777777777778
$ws.Range("J2").Value = 0.02777777777777778
$ws.Range("P2").Value = 0.1215277777777778
$ws.Range("S2").Value = 0.08680555555555555

# Row 3
$ws.Range("B3").Value = 0.005617977528089887
$ws.Range("C3").Value = 0.02247191011235955
$ws.Range("J3").Value = 0.02808988764044944
$ws.Range("P3").Value = 0.7471910112359551
$ws.Range("S3").Value = 0.1966292134831461

# Row 4
$ws.Range("J4").Value = 0.07317073170731707
$ws.Range("P4").Value = 0.6585365853658537
$ws.Range("S4").Value = 0.2682926829268293

# Row 5
$ws.Range("P5").Value = 1

# Row 6
$ws.Range("B6").Value = 0.08181818181818182
$ws.Range("D6").Value = 0.004545454545454545
$ws.Range("F6").Value = 0.06363636363636363
$ws.Range("J6").Value = 0.2318181818181818
$ws.Range("O6").Value = 0.01363636363636364
$ws.Range("Q6").Value = 0.1772727272727273
$ws.Range("R6").Value = 0.08636363636363636
$ws.Range("S6").Value = 0.3409090909090909

# Row 7
$ws.Range("B7").Value = 0.1602209944751381
$ws.Range("D7").Value = 0.02209944751381215
$ws.Range("F7").Value = 0.05524861878453038
$ws.Range("J7").Value = 0.1104972375690608
$ws.Range("O7").Value = 0.01104972375690608
$ws.Range("Q7").Value = 0.1491712707182321
$ws.Range("R7").Value = 0.09392265193370165
$ws.Range("S7").Value = 0.3977900552486188

# Row 8
$ws.Range("B8").Value = 0.0968586387434555
$ws.Range("D8").Value = 0.02094240837696335
$ws.Range("F8").Value = 0.06282722513089005
$ws.Range("J8").Value = 0.1282722513089005
$ws.Range("O8").Value = 0.01832460732984293
$ws.Range("Q8").Value = 0.1884816753926702
$ws.Range("R8").Value = 0.09947643979057591
$ws.Range("S8").Value = 0.3848167539267016

# Row 9
$ws.Range("B9").Value = 0.09793814432989691
$ws.Range("D9").Value = 0.0154639175257732
$ws.Range("E9").Value = 0.005154639175257732
$ws.Range("F9").Value = 0.06701030927835051
$ws.Range("J9").Value = 0.1134020618556701
$ws.Range("O9").Value = 0.01030927835051546
$ws.Range("Q9").Value = 0.2010309278350516
$ws.Range("R9").Value = 0.06185567010309279
$ws.Range("S9").Value = 0.4278350515463917

# Row 10
$ws.Range("B10").Value = 0.1076417419884963
$ws.Range("D10").Value = 0.02136400986031224
$ws.Range("F10").Value = 0.08463434675431389
$ws.Range("J10").Value = 0.1035332785538209
$ws.Range("O10").Value = 0.01314708299096138
$ws.Range("Q10").Value = 0.1914543960558751
$ws.Range("R10").Value = 0.08299096138044372
$ws.Range("S10").Value = 0.3952341824157765

# Row 11
$ws.Range("G11").Value = 0.1701388888888889
$ws.Range("J11").Value = 0.0763888888888889
$ws.Range("K11").Value = 0.2326388888888889
$ws.Range("L11").Value = 0.5
$ws.Range("S11").Value = 0.02083333333333333

# Row 12
$ws.Range("G12").Value = 0.7379310344827587
$ws.Range("J12").Value = 0.2068965517241379
$ws.Range("K12").Value = 0.01379310344827586
$ws.Range("L12").Value = 0.01379310344827586
$ws.Range("S12").Value = 0.02758620689655172

# Row 13
$ws.Range("G13").Value = 0.6190476190476191
$ws.Range("J13").Value = 0.3571428571428572
$ws.Range("S13").Value = 0.02380952380952381

# Row 15
$ws.Range("F15").Value = 0.02105263157894737
$ws.Range("H15").Value = 0.1421052631578947
$ws.Range("I15").Value = 0.09473684210526316
$ws.Range("J15").Value = 0.4052631578947368
$ws.Range("K15").Value = 0.03684210526315789
$ws.Range("M15").Value = 0.01578947368421053
$ws.Range("O15").Value = 0.04210526315789474
$ws.Range("S15").Value = 0.2421052631578947

# Row 16
$ws.Range("F16").Value = 0.005263157894736842
$ws.Range("H16").Value = 0.2052631578947368
$ws.Range("I16").Value = 0.09473684210526316
$ws.Range("J16").Value = 0.4052631578947368
$ws.Range("K16").Value = 0.08421052631578947
$ws.Range("M16").Value = 0.01578947368421053
$ws.Range("O16").Value = 0.03684210526315789
$ws.Range("S16").Value = 0.1526315789473684

# Row 17
$ws.Range("F17").Value = 0.019559902200489
$ws.Range("H17").Value = 0.1613691931540342
$ws.Range("I17").Value = 0.1026894865525672
$ws.Range("J17").Value = 0.4449877750611247
$ws.Range("K17").Value = 0.08557457212713937
$ws.Range("M17").Value = 0.01222493887530562
$ws.Range("O17").Value = 0.06112469437652811
$ws.Range("S17").Value = 0.1124694376528117

# Row 18
$ws.Range("F18").Value = 0.0106951871657754
$ws.Range("H18").Value = 0.1711229946524064
$ws.Range("I18").Value = 0.106951871657754
$ws.Range("J18").Value = 0.4705882352941176
$ws.Range("K18").Value = 0.1122994652406417
$ws.Range("M18").Value = 0.0160427807486631
$ws.Range("N18").Value = 0.0053475935828877
$ws.Range("O18").Value = 0.0427807486631016
$ws.Range("S18").Value = 0.06417112299465241

# Row 19
$ws.Range("F19").Value = 0.01261564339781329
$ws.Range("H19").Value = 0.1867115222876367
$ws.Range("I19").Value = 0.08242220353238015
$ws.Range("J19").Value = 0.3851976450798991
$ws.Range("K19").Value = 0.1177460050462574
$ws.Range("M19").Value = 0.02523128679562658
$ws.Range("N19").Value = 0.002523128679562658
$ws.Range("O19").Value = 0.07569386038687972
